# Append new row 29 to the "Artfynd" sheet, matching the data in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Plain numeric cells
$ws.Cells.Item($row, 1).Value = 111858253            # A - Id
$ws.Cells.Item($row, 2).Value = 89802                # B - Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value = 5420                 # E - TaxonId
$ws.Cells.Item($row, 17).Value = 682892.0984107114   # Q - Ost
$ws.Cells.Item($row, 18).Value = 6575487.713710153   # R - Nord
$ws.Cells.Item($row, 19).Value = 50                  # S - Noggrannhet

# Plain text cells (non numeric-looking, safe to assign directly)
$ws.Cells.Item($row, 3).Value = "Ovaliderad"                    # C - Valideringsstatus
$ws.Cells.Item($row, 4).Value = "LC"                             # D - Rodlistade
$ws.Cells.Item($row, 6).Value = "Grovticka"                      # F - Artnamn
$ws.Cells.Item($row, 7).Value = "Phaeolus schweinitzii"          # G - Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Fr.) Pat."                     # H - Auktor
$ws.Cells.Item($row, 10).Value = "fruktkroppar"                  # J - Enhet
$ws.Cells.Item($row, 16).Value = "NW Erstavik, Srm"               # P - Lokalnamn
$ws.Cells.Item($row, 20).Value = "Stockholm"                     # T - Lan
$ws.Cells.Item($row, 21).Value = "Nacka"                          # U - Kommun
$ws.Cells.Item($row, 22).Value = "Södermanland"                  # V - Provins
$ws.Cells.Item($row, 23).Value = "Nacka"                          # W - Forsamling
$ws.Cells.Item($row, 29).Value = 'Hällmark, stig; På Tallstubbe, "tung-formad"'  # AC - Publik kommentar
$ws.Cells.Item($row, 49).Value = "Henry Gudmundson"               # AW - Rapportor
$ws.Cells.Item($row, 50).Value = "Henry Gudmundson"               # AX - Observatorer

# Text cells that look like numbers/dates - force text format so Excel
# does not silently reinterpret them as numeric or date values.
$textForced = @(9, 25, 26, 27, 28)
foreach ($c in $textForced) {
    $ws.Cells.Item($row, $c).NumberFormat = "@"
}
$ws.Cells.Item($row, 9).Value = "1"                    # I - Antal
$ws.Cells.Item($row, 25).Value = "2023-09-02"          # Y - Startdatum
$ws.Cells.Item($row, 26).Value = "00:00"               # Z - Starttid
$ws.Cells.Item($row, 27).Value = "2023-09-02"          # AA - Slutdatum
$ws.Cells.Item($row, 28).Value = "00:00"               # AB - Sluttid

# Boolean cells
$ws.Cells.Item($row, 30).Value = $false   # AD - Ej aterfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE - Osaker artbestamning
$ws.Cells.Item($row, 33).Value = $false   # AG - Ospontan
